$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: extend header with two new columns (P1, Q1) ---
# Copy formatting (bold/border/center) from O1 to P1:Q1
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Rows 2-25: update existing columns I, K, M, O and add new columns P, Q ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
